$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates that apply identically to the
# "展览" sheet (sheet1) and the "全部类型" sheet (sheet4).
$updates = @{
    2  = 1567
    3  = 48
    4  = 1023
    5  = 23
    7  = 2612
    9  = 1651
    12 = 542
    15 = 59
    16 = 78
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
